# Update TPM-derived statistics in the LR-pairs sheet (Fn1-Plaur)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 2.843949
$ws.Range("N2").Value = 8.531846999999999
$ws.Range("O2").Value = 0.4976240243095911
$ws.Range("P2").Value = 0.4976240243095912
$ws.Range("Q2").Value = 17.746548906492
$ws.Range("R2").Value = 159.718940158428
$ws.Range("S2").Value = 0.008619995241555741
$ws.Range("T2").Value = 0.008619995241555744

# Row 3
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.4403664892852895
$ws.Range("P3").Value = 0.4403664892852897
$ws.Range("Q3").Value = 15.704598365652
$ws.Range("R3").Value = 141.341385290868
$ws.Range("S3").Value = 0.007628162742838541
$ws.Range("T3").Value = 0.007628162742838544

# Row 4
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("M4").Value = 0.3543876666666666
$ws.Range("O4").Value = 0.06200948640511928
$ws.Range("P4").Value = 0.0620094864051193
$ws.Range("Q4").Value = 2.211417313868
$ws.Range("R4").Value = 19.90275582481199
$ws.Range("S4").Value = 0.001074147250999476
$ws.Range("T4").Value = 0.001074147250999477

# Row 5
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("M5").Value = 2.843949
$ws.Range("N5").Value = 8.531846999999999
$ws.Range("O5").Value = 0.4976240243095911
$ws.Range("P5").Value = 0.4976240243095912
$ws.Range("Q5").Value = 982.773726780471
$ws.Range("R5").Value = 8844.963541024237
$ws.Range("S5").Value = 0.4773606909721268
$ws.Range("T5").Value = 0.477360690972127

# Row 6
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.4403664892852895
$ws.Range("P6").Value = 0.4403664892852897
$ws.Range("S6").Value = 0.422434692331923
$ws.Range("T6").Value = 0.4224346923319232

# Row 7
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("M7").Value = 0.3543876666666666
$ws.Range("O7").Value = 0.06200948640511928
$ws.Range("P7").Value = 0.0620094864051193
$ws.Range("S7").Value = 0.05948444976755903
$ws.Range("T7").Value = 0.05948444976755906

# Row 8
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("M8").Value = 2.843949
$ws.Range("N8").Value = 8.531846999999999
$ws.Range("O8").Value = 0.4976240243095911
$ws.Range("P8").Value = 0.4976240243095912
$ws.Range("Q8").Value = 23.970902902328
$ws.Range("R8").Value = 215.738126120952
$ws.Range("S8").Value = 0.01164333809590852
$ws.Range("T8").Value = 0.01164333809590853

# Row 9
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.4403664892852895
$ws.Range("P9").Value = 0.4403664892852897
$ws.Range("S9").Value = 0.01030363421052797
$ws.Range("T9").Value = 0.01030363421052797

# Row 10
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("M10").Value = 0.3543876666666666
$ws.Range("O10").Value = 0.06200948640511928
$ws.Range("P10").Value = 0.0620094864051193
$ws.Range("Q10").Value = 2.98704102902311
$ws.Range("R10").Value = 26.883369261208
$ws.Range("S10").Value = 0.001450889386560775
$ws.Range("T10").Value = 0.001450889386560776
